$d = $word.ActiveDocument

# --- Step 1: fix "Descibe" run split (merge "Descibe" + " " into "Descibe ") ---
$rng = $d.Content
$rng.Find.Execute("Descibe", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$descEnd = $rng.End
$rng.InsertAfter(" ")
# delete the now-redundant standalone space run that followed "Descibe"
$dupSpace = $d.Range($descEnd + 1, $descEnd + 2)
$dupSpace.Text = ""

# --- Step 2: round-trip the whole document content through OpenXML.
#     This merges runs that share identical formatting (e.g. "for training/" + "voc"
#     + " counseling, ...") and drops now-stale proofErr spell-check markers. ---
$full = $d.Content
$full.InsertXML($full.WordOpenXML)

# --- Step 3: insert a new, completely empty paragraph right after the main
#     paragraph (i.e. right before the paragraph that holds the _GoBack bookmark) ---
$bookmarkPara = $d.Paragraphs.Item(2)
$bookmarkPara.Range.InsertBefore("`r")

# round-trip again so the freshly split paragraph collapses down to a bare <w:p/>
$full = $d.Content
$full.InsertXML($full.WordOpenXML)

# --- Step 4: add a "4 pages" run at the start of the bookmark paragraph (now #3) ---
$bookmarkPara = $d.Paragraphs.Item(3)
$bookmarkPara.Range.InsertBefore("4 pages")
